{"js": "// Replace the 25 division-problem answers in the table with new values.\n// Each old text is unique within the document, so an exact, case-sensitive\n// search-and-replace on each pair is sufficient and preserves all run\n// formatting (font, size, etc.) since only the text inside the existing\n// run is swapped out.\nconst pairs = [\n  [\"817\u00f75=163, 2\", \"561\u00f78=70, 1\"],\n  [\"772\u00f75=154, 2\", \"676\u00f74=169, 0\"],\n  [\"353\u00f76=58, 5\", \"751\u00f76=125, 1\"],\n  [\"303\u00f79=33, 6\", \"933\u00f79=103, 6\"],\n  [\"530\u00f73=176, 2\", \"152\u00f72=76, 0\"],\n  [\"854\u00f73=284, 2\", \"753\u00f72=376, 1\"],\n  [\"576\u00f78=72, 0\", \"710\u00f76=118, 2\"],\n  [\"584\u00f72=292, 0\", \"958\u00f77=136, 6\"],\n  [\"975\u00f74=243, 3\", \"781\u00f73=260, 1\"],\n  [\"312\u00f78=39, 0\", \"540\u00f72=270, 0\"],\n  [\"368\u00f77=52, 4\", \"495\u00f77=70, 5\"],\n  [\"987\u00f73=329, 0\", \"836\u00f79=92, 8\"],\n  [\"645\u00f78=80, 5\", \"557\u00f78=69, 5\"],\n  [\"235\u00f78=29, 3\", \"285\u00f74=71, 1\"],\n  [\"679\u00f74=169, 3\", \"266\u00f79=29, 5\"],\n  [\"914\u00f79=101, 5\", \"633\u00f75=126, 3\"],\n  [\"705\u00f75=141, 0\", \"182\u00f79=20, 2\"],\n  [\"196\u00f72=98, 0\", \"584\u00f79=64, 8\"],\n  [\"297\u00f77=42, 3\", \"461\u00f78=57, 5\"],\n  [\"602\u00f74=150, 2\", \"485\u00f76=80, 5\"],\n  [\"921\u00f79=102, 3\", \"518\u00f74=129, 2\"],\n  [\"348\u00f79=38, 6\", \"553\u00f77=79, 0\"],\n  [\"462\u00f72=231, 0\", \"627\u00f73=209, 0\"],\n  [\"204\u00f75=40, 4\", \"371\u00f72=185, 1\"],\n  [\"245\u00f72=122, 1\", \"472\u00f79=52, 4\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answers in the table with new values.\n# Each \"old\" text is unique within the document, so a straightforward\n# Find/Replace (restricted to exact text, case-sensitive, no wildcards)\n# for each pair is sufficient and preserves all existing run formatting\n# (font, size, etc.), because Word's Find.Execute replaces only the\n# matched text run-locally.\n\n$d = $word.ActiveDocument\n$c = [char]0x00F7   # '\u00f7' (DIVISION SIGN, U+00F7) - built via interpolation\n                     # below to dodge PowerShell's \"numeric-looking string\"\n                     # + [char] => arithmetic-addition coercion quirk.\n\n$pairs = @(\n  ,@(\"817${c}5=163, 2\", \"561${c}8=70, 1\")\n  ,@(\"772${c}5=154, 2\", \"676${c}4=169, 0\")\n  ,@(\"353${c}6=58, 5\", \"751${c}6=125, 1\")\n  ,@(\"303${c}9=33, 6\", \"933${c}9=103, 6\")\n  ,@(\"530${c}3=176, 2\", \"152${c}2=76, 0\")\n  ,@(\"854${c}3=284, 2\", \"753${c}2=376, 1\")\n  ,@(\"576${c}8=72, 0\", \"710${c}6=118, 2\")\n  ,@(\"584${c}2=292, 0\", \"958${c}7=136, 6\")\n  ,@(\"975${c}4=243, 3\", \"781${c}3=260, 1\")\n  ,@(\"312${c}8=39, 0\", \"540${c}2=270, 0\")\n  ,@(\"368${c}7=52, 4\", \"495${c}7=70, 5\")\n  ,@(\"987${c}3=329, 0\", \"836${c}9=92, 8\")\n  ,@(\"645${c}8=80, 5\", \"557${c}8=69, 5\")\n  ,@(\"235${c}8=29, 3\", \"285${c}4=71, 1\")\n  ,@(\"679${c}4=169, 3\", \"266${c}9=29, 5\")\n  ,@(\"914${c}9=101, 5\", \"633${c}5=126, 3\")\n  ,@(\"705${c}5=141, 0\", \"182${c}9=20, 2\")\n  ,@(\"196${c}2=98, 0\", \"584${c}9=64, 8\")\n  ,@(\"297${c}7=42, 3\", \"461${c}8=57, 5\")\n  ,@(\"602${c}4=150, 2\", \"485${c}6=80, 5\")\n  ,@(\"921${c}9=102, 3\", \"518${c}4=129, 2\")\n  ,@(\"348${c}9=38, 6\", \"553${c}7=79, 0\")\n  ,@(\"462${c}2=231, 0\", \"627${c}3=209, 0\")\n  ,@(\"204${c}5=40, 4\", \"371${c}2=185, 1\")\n  ,@(\"245${c}2=122, 1\", \"472${c}9=52, 4\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    $result = $find.Execute(\n        [ref]$oldText,    # FindText\n        [ref]$false,      # MatchCase\n        [ref]$false,      # MatchWholeWord\n        [ref]$false,      # MatchWildcards\n        [ref]$false,      # MatchSoundsLike\n        [ref]$false,      # MatchAllWordForms\n        [ref]$true,       # Forward\n        1,                # Wrap (wdFindContinue)\n        [ref]$false,      # Format\n        [ref]$newText,    # ReplaceWith\n        2                 # Replace (wdReplaceAll)\n    )\n\n    if (-not $result) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
